$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# 1) Rename header labels: "<name>_old" -> "<name>_FV2410", "<name>_new" -> "<name>_FV2504"
#    (AHB comparison columns were generated against FV2410 vs FV2504 message versions)
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $label = $cell.Text
    if ($label -like "*_old") {
        $cell.Value = ($label -replace "_old$", "_FV2410")
    } elseif ($label -like "*_new") {
        $cell.Value = ($label -replace "_new$", "_FV2504")
    }
}

# 2) Turn the used range into an Excel Table (ListObject) with a header row
$tbl = $ws.ListObjects.Add(1, $used, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze the header row (split/freeze at row 2)
$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true
